# FOIA all-cause / covid deaths workbook — add quarterly "CY" total rows
# and leave the selection/active-sheet state where the editor left it
# (covid deaths tab active, row-13 totals just entered and still selected).

$wb = $excel.ActiveWorkbook

# --- "all-cause deaths" sheet: totals across every age-group row (3-12) ---
$ws1 = $wb.Worksheets.Item("all-cause deaths")
$ws1.Range("F13").Formula  = "=SUM(F3:F12)"
$ws1.Range("K13").Formula  = "=SUM(K3:K12)"
$ws1.Range("P13").Formula  = "=SUM(P3:P12)"
$ws1.Range("U13").Formula  = "=SUM(U3:U12)"
$ws1.Range("Z13").Formula  = "=SUM(Z3:Z12)"
$ws1.Range("AE13").Formula = "=SUM(AE3:AE12)"

# --- "covid deaths" sheet: totals starting at row 5 (no COVID deaths before) ---
$ws2 = $wb.Worksheets.Item("covid deaths")
$ws2.Range("F13").Formula = "=SUM(F5:F12)"
$ws2.Range("K13").Formula = "=SUM(K5:K12)"
$ws2.Range("P13").Formula = "=SUM(P5:P12)"
$ws2.Range("U13").Formula = "=SUM(U5:U12)"

# Mirror the author's final selection state: the all-cause sheet was edited
# first (selection left on its new total row), then the covid sheet was
# edited and left active/selected.
$ws1.Activate()
$ws1.Range("F13:AE13").Select()

$ws2.Activate()
$ws2.Range("F13:U13").Select()
